$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# Step 1: capture the 6 existing rows (2..7) of "New" - values + hyperlink
# URLs - before any mutation happens, while rId1..rId6 still line up with
# rows 2..7.
# ---------------------------------------------------------------------------
$moveCount = 6
$srcVals = @()
$srcUrls = @()
for ($i = 0; $i -lt $moveCount; $i++) {
    $r = 2 + $i
    $row = @(
        $wsNew.Range("A$r").Value2,
        $wsNew.Range("B$r").Value2,
        $wsNew.Range("C$r").Value2,
        $wsNew.Range("D$r").Value2,
        $wsNew.Range("E$r").Value2,
        $wsNew.Range("F$r").Value2
    )
    $srcVals += ,$row
    $srcUrls += $wsNew.Hyperlinks.Item($i + 1).Address
}

# ---------------------------------------------------------------------------
# Step 2: append those 6 rows onto the end of "Previously added" (rows
# 323..328), copying the row-322 formatting down so the cell styles match
# the rest of the sheet.
# ---------------------------------------------------------------------------
$destStart = 323
for ($i = 0; $i -lt $moveCount; $i++) {
    $r = $destStart + $i
    $row = $srcVals[$i]

    $wsPrev.Range("A322:F322").Copy($wsPrev.Range("A${r}:F${r}"))

    $wsPrev.Range("A$r").Value = $row[0]
    $wsPrev.Range("B$r").Value = $row[1]
    $wsPrev.Range("C$r").Value = $row[2]
    $wsPrev.Range("D$r").Value = $row[3]
    $wsPrev.Range("E$r").Value = $row[4]
    $wsPrev.Range("F$r").Value = $row[5]

    $wsPrev.Hyperlinks.Add($wsPrev.Range("A$r"), $srcUrls[$i])
    # Hyperlinks.Add forces Excel's built-in "Hyperlink" style onto the
    # cell; restore the sheet's own link-cell style (same as row 322) by
    # re-copying the format from a known-good cell, then re-apply the value
    # (copy also clobbers the cell contents).
    $wsPrev.Range("A322").Copy($wsPrev.Range("A$r"))
    $wsPrev.Range("A$r").Value = $row[0]
}

# ---------------------------------------------------------------------------
# Step 3: "New" loses its current 6 listings and gains 2 brand-new ones.
# Wipe every hyperlink on the sheet (rows 2..7 all disappear/move), delete
# the now-unneeded rows 4..7, then populate rows 2..3 with the new data.
# ---------------------------------------------------------------------------
$wsNew.Range("A2").Hyperlinks.Delete()
$wsNew.Rows("4:7").Delete()

$newRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/kcjhx.html",
      "68 000 €", "Cēsis un raj.", "17 ha.", "42760020065", 46000.52291666667),
    @("https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/makonkalna-pag/gohej.html",
      "190 000 €", "Rēzekne un raj.", "20 ha.", "78720020066", 45999.61041666666)
)

for ($i = 0; $i -lt 2; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]

    $wsNew.Range("A$r").Value = $row[0]
    $wsNew.Range("B$r").Value = $row[1]
    $wsNew.Range("C$r").Value = $row[2]
    $wsNew.Range("D$r").Value = $row[3]
    $wsNew.Range("E$r").Value = $row[4]
    $wsNew.Range("F$r").Value = $row[5]

    $wsNew.Hyperlinks.Add($wsNew.Range("A$r"), $row[0])
    # restore the original link-cell style (style index 3), borrowed from
    # "Previously added" which still has plenty of untouched style-3 cells.
    $wsPrev.Range("A2").Copy($wsNew.Range("A$r"))
    $wsNew.Range("A$r").Value = $row[0]
}
